# automatic bills dowloader 1.1
#
# Appends the newly-downloaded consumer/unit rows to Sheet1 and gives the
# new block its "billed" look: a light-grey box border, right aligned +
# wrapped text, and a custom 0/(0) number format on the consumer numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- data for the new rows (consumer number, unit) --------------------
$rows = @(
    @(160252832341, 4610),
    @(160252832367, 4610),
    @(170743925516, 4746),
    @(170741644315, 4746),
    @(170741644412, 4746),
    @(170741644404, 4746),
    @(170003628933, 4745),
    @(170003629018, 4745),
    @(170003629026, 4745),
    @(170003629034, 4745)
)

$firstRow = 16
$lastRow = $firstRow + $rows.Length - 1

# Row 15 is a blank spacer row, same height as the data rows.
$ws.Rows(15).RowHeight = 15.75

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Rows($r).RowHeight = 15.75
}

# Start from the same font the existing data rows use (Arial 10, theme text
# colour) before layering the new box-border / alignment on top of it.
$ws.Range("A2").Copy()
$ws.Range("A16:B25").PasteSpecial(-4122)

# -- formatting for the new block --------------------------------------
$block = $ws.Range("A16:B25")
$block.Borders.Weight = -4138
$block.Borders.Color = 13421772
$block.HorizontalAlignment = -4152
$block.WrapText = $true

$ws.Range("A16:A25").NumberFormat = "0_);[Red]\(0\)"

# -- selection tweak, matches where the user left off -------------------
[void]$ws.Range("D13").Select()

Write-Output "done"
